$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42606.88181712963
$ws.Range("B5").Value = 36
$ws.Range("C5").Value = 71
$ws.Range("D5").Value = 27
$ws.Range("E5").Value = 58
$ws.Range("F5").Value = 41
$ws.Range("G5").Value = 8246
$ws.Range("H5").Value = 8918
$ws.Range("I5").Value = 1678
$ws.Range("J5").Value = 222
$ws.Range("K5").Value = 86
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = "Noun"
